# add deep learning psets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in filenames first (column A), top to bottom
$ws.Cells.Item(64, 1).Value = "dl-ps1"
$ws.Cells.Item(65, 1).Value = "dl-ps2"
$ws.Cells.Item(66, 1).Value = "dl-ps3"
$ws.Cells.Item(67, 1).Value = "dl-ps4"

# Fill in titles (column B)
$ws.Cells.Item(64, 2).Value = "Training a MLP"
$ws.Cells.Item(65, 2).Value = "Training a CNN"
$ws.Cells.Item(67, 2).Value = "Training a GAN/VAE"
$ws.Cells.Item(66, 2).Value = "Thermodynamics of SGD"

# Fill in subjects (column C) - reuse existing "deep-learning" string
$ws.Cells.Item(64, 3).Value = "deep-learning"
$ws.Cells.Item(65, 3).Value = "deep-learning"
$ws.Cells.Item(66, 3).Value = "deep-learning"
$ws.Cells.Item(67, 3).Value = "deep-learning"

# Fill in groups (column D) - reuse existing "CMSC-31230" string
$ws.Cells.Item(64, 4).Value = "CMSC-31230"
$ws.Cells.Item(65, 4).Value = "CMSC-31230"
$ws.Cells.Item(66, 4).Value = "CMSC-31230"
$ws.Cells.Item(67, 4).Value = "CMSC-31230"

# Match the styling used on the other CMSC-31230 group rows (explicit black font)
$ws.Cells.Item(65, 4).Font.Color = 0
$ws.Cells.Item(66, 4).Font.Color = 0
$ws.Cells.Item(67, 4).Font.Color = 0

# Update the view to reflect where the user ended up after the edit
$ws.Application.ActiveWindow.ScrollRow = 45
$ws.Range("B66").Select()
